$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Historias de usuario")

# Update wording of existing user stories. Order matters for shared-string
# table insertion order, so follow the same sequence the original author
# used: D7, D8, D5, then D6/D9/E11 (which reuse already-existing strings),
# then the new row's cells C12, D12, E12, B12.
$ws.Range("D7").Value = "Necesito saber el motivo por el cual no puedo acceder a mi cuenta de usuario debido a inactividad, mal comportamiento u otros casos; además, me permita proceder a una solución si hubiera oportunidad."
$ws.Range("D8").Value = "Necesito eliminar cuentas de usuarios inactivos que no cuentan con una membresia de pago."
$ws.Range("D5").Value = "Necesito ver un panel de control con un menu de opciones; y pueda ver mi perfil de usuario, mis mascotas y configuraciones."
$ws.Range("D6").Value = "Necesito realizar busquedas de canes según los filtros que establezca para la busqueda que pueden ser: raza, sexo, edad, peso y  ubicación del dueño."
$ws.Range("D9").Value = "Necesito bloquear usuarios que presenten un comportamiento indebido para evitar la mala experiencia de los usuarios en el sitio. Dentro de la web se considera mal comportamiento. Por ejemplo, escribir una publicación ofendiendo a la mascota o al dueño. Proporcionar información falsa en alguna publicación realizada. Denunciar de forma errónea y/o sin responsabilidad algúna al usuario sin que este allá cometido algún mal comportamiento."
$ws.Range("E11").Value = "Con la finalidad de asegurar al usuario si debe proceder a un acuerdo con el dueño de la mascota candidata."

# Add a new row (12) for the 2nd sprint's user story (HU10)
$ws.Range("C12").Value = "Como un usuario o administrador"
$ws.Range("D12").Value = "Necesito poder cerrar la sesion de mi cuenta logeada"
$ws.Range("E12").Value = "Con la finalidad de finalizar mis operaciónes y asegurar el cierre de mi sesion en la web."
$ws.Range("B12").Value = "HU10"

# Copy formatting from the equivalent cells in the existing table so the
# new row matches the sheet's styling (borders/alignment/wrap).
$ws.Range("B11").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("C11").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("E11").Copy()
$ws.Range("E12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows.Item(12).RowHeight = 45

# Update sheet view selection/scroll like the final saved state
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("D16").Select()
